$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price text in column D is preserved as text
# (avoids Excel auto-converting strings like '141.40' or '0.00001028' into numbers)
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

# Apply the updated coin data
$ws.Cells.Item(2, 4).Value = "20.525.05"
$ws.Cells.Item(2, 5).Value = "  +2.10%  "
$ws.Cells.Item(3, 4).Value = "1.474.32"
$ws.Cells.Item(3, 5).Value = "  +3.65%  "
$ws.Cells.Item(4, 4).Value = "1.008"
$ws.Cells.Item(4, 5).Value = "  +0.78%  "
$ws.Cells.Item(5, 4).Value = "0.9571"
$ws.Cells.Item(5, 5).Value = "  -4.00%  "
$ws.Cells.Item(6, 4).Value = "276.76"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "
$ws.Cells.Item(7, 4).Value = "0.3649"
$ws.Cells.Item(7, 5).Value = "  -1.60%  "
$ws.Cells.Item(8, 5).Value = "  -2.89%  "
$ws.Cells.Item(9, 4).Value = "39.72"
$ws.Cells.Item(9, 5).Value = "  +0.58%  "
$ws.Cells.Item(10, 4).Value = "1.055"
$ws.Cells.Item(10, 5).Value = "  -0.56%  "
$ws.Cells.Item(11, 4).Value = "0.06613"
$ws.Cells.Item(11, 5).Value = "  +0.91%  "
$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  +0.33%  "
$ws.Cells.Item(13, 4).Value = "18.15"
$ws.Cells.Item(13, 5).Value = "  +1.41%  "
$ws.Cells.Item(14, 4).Value = "5.464"
$ws.Cells.Item(14, 5).Value = "  -1.17%  "
$ws.Cells.Item(15, 5).Value = "  -0.43%  "
$ws.Cells.Item(16, 4).Value = "0.00001028"
$ws.Cells.Item(16, 5).Value = "  +0.51%  "
$ws.Cells.Item(17, 4).Value = "1.474.51"
$ws.Cells.Item(17, 5).Value = "  +3.67%  "
$ws.Cells.Item(18, 4).Value = "0.05895"
$ws.Cells.Item(18, 5).Value = "  +3.11%  "
$ws.Cells.Item(19, 4).Value = "0.9647"
$ws.Cells.Item(19, 5).Value = "  -3.26%  "
$ws.Cells.Item(20, 4).Value = "69.27"
$ws.Cells.Item(20, 5).Value = "  -3.40%  "
$ws.Cells.Item(21, 4).Value = "5.468"
$ws.Cells.Item(21, 5).Value = "  -2.66%  "
$ws.Cells.Item(22, 4).Value = "14.48"
$ws.Cells.Item(22, 5).Value = "  -2.65%  "
$ws.Cells.Item(23, 4).Value = "11.02"
$ws.Cells.Item(23, 5).Value = "  -0.56%  "
$ws.Cells.Item(24, 4).Value = "2.252"
$ws.Cells.Item(24, 5).Value = "  +1.20%  "
$ws.Cells.Item(25, 4).Value = "20.584.21"
$ws.Cells.Item(25, 5).Value = "  +2.23%  "
$ws.Cells.Item(26, 4).Value = "141.40"
$ws.Cells.Item(26, 5).Value = "  +5.09%  "
$ws.Cells.Item(27, 4).Value = "2.131"
$ws.Cells.Item(27, 5).Value = "  -7.05%  "
$ws.Cells.Item(28, 4).Value = "17.20"
$ws.Cells.Item(28, 5).Value = "  -0.64%  "
$ws.Cells.Item(29, 4).Value = "1.632.24"
$ws.Cells.Item(29, 5).Value = "  +3.19%  "
$ws.Cells.Item(30, 4).Value = "113.52"
$ws.Cells.Item(30, 5).Value = "  +2.28%  "
$ws.Cells.Item(31, 4).Value = "3.940"
$ws.Cells.Item(31, 5).Value = "  -0.60%  "
$ws.Cells.Item(32, 4).Value = "0.8204"
$ws.Cells.Item(32, 5).Value = "  -1.25%  "
$ws.Cells.Item(33, 4).Value = "4.974"
$ws.Cells.Item(33, 5).Value = "  -5.95%  "
$ws.Cells.Item(34, 4).Value = "0.07949"
$ws.Cells.Item(34, 5).Value = "  +1.88%  "
$ws.Cells.Item(35, 4).Value = "1.532"
$ws.Cells.Item(35, 5).Value = "  +3.62%  "
$ws.Cells.Item(36, 4).Value = "1.235"
$ws.Cells.Item(36, 5).Value = "  +11.26%  "
$ws.Cells.Item(37, 4).Value = "0.05759"
$ws.Cells.Item(37, 5).Value = "  -1.93%  "
$ws.Cells.Item(38, 4).Value = "4.725"
$ws.Cells.Item(38, 5).Value = "  -4.09%  "
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "0.02038"
$ws.Cells.Item(39, 5).Value = "  -1.41%  "
$ws.Cells.Item(40, 2).Value = "Aptos"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(40, 4).Value = "10.44"
$ws.Cells.Item(40, 5).Value = "  -1.62%  "
$ws.Cells.Item(41, 2).Value = "Frax"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(41, 4).Value = "0.9573"
$ws.Cells.Item(41, 5).Value = "  -3.95%  "
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "7.589"
$ws.Cells.Item(42, 5).Value = "  -5.37%  "
$ws.Cells.Item(43, 4).Value = "0.1879"
$ws.Cells.Item(43, 5).Value = "  +0.16%  "
$ws.Cells.Item(44, 4).Value = "0.5287"
$ws.Cells.Item(44, 5).Value = "  -1.15%  "
$ws.Cells.Item(45, 4).Value = "3.506"
$ws.Cells.Item(45, 5).Value = "  -1.29%  "
$ws.Cells.Item(46, 4).Value = "12.11"
$ws.Cells.Item(46, 5).Value = "  -1.97%  "
$ws.Cells.Item(47, 4).Value = "117.37"
$ws.Cells.Item(47, 5).Value = "  -0.87%  "
$ws.Cells.Item(48, 4).Value = "0.5188"
$ws.Cells.Item(48, 5).Value = "  -0.99%  "
$ws.Cells.Item(49, 4).Value = "1.776"
$ws.Cells.Item(49, 5).Value = "  -0.47%  "
$ws.Cells.Item(50, 4).Value = "0.06456"
$ws.Cells.Item(50, 5).Value = "  +3.48%  "
$ws.Cells.Item(51, 4).Value = "0.9935"
$ws.Cells.Item(51, 5).Value = "  -0.35%  "
